# Apply the 2017-sheet "dividend received" update.
# Only the raw input cells (share price in column D, and the dividend
# amount landing in column M for a few rows) are written; every other
# touched cell (E, T, U, V and the row-14 totals) is a formula that
# recalculates on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")
$ws.Activate()

# Row 2
$ws.Range("D2").Value = 22.254000000000001
$ws.Range("M2").Value = 4.21

# Row 4
$ws.Range("D4").Value = 20.681999999999999
$ws.Range("M4").Value = 9.85

# Row 10
$ws.Range("D10").Value = 8.1370000000000005
$ws.Range("M10").Value = 6.07

# Row 12
$ws.Range("D12").Value = 15.423999999999999
$ws.Range("M12").Value = 7.82

# Row 13
$ws.Range("D13").Value = 12.224
$ws.Range("M13").Value = 9.32

# Make sure everything is up to date before the selection/save.
$excel.CalculateFull()

# Update the selection shown when the sheet was last saved.
$ws.Range("F19").Select()
